$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.120.67'
$ws.Range('E2').Value = '  +2.21%  '
$ws.Range('D3').Value = '3.775.02'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.97%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '625.25'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.49'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.28%  '
$ws.Range('D7').Value = '3.772.14'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +1.60%  '
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.457'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.70'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.68'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.87%  '
$ws.Range('D15').Value = '4.414.20'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '3.778.22'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('D17').Value = '69.160.89'
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.64'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.10'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.83%  '
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '468.31'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.62'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.95%  '
$ws.Range('E23').Value = '  +2.57%  '
$ws.Range('E24').Value = '  +4.32%  '
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.04'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.93%  '
$ws.Range('E27').Value = '  +3.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.02'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.07%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '3.925.41'
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('E31').Value = '  +3.75%  '
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.23'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.75'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.727.61'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.99'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.165'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +15.43%  '
$ws.Range('E39').Value = '  +2.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.42'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +8.53%  '
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.968'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.297'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.05'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '152.73'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('E48').Value = '  +4.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '46.65'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('E50').Value = '  +1.92%  '
$ws.Range('E51').Value = '  -0.47%  '
